$d = $word.ActiveDocument

# --- Step 1: "Debug using KF" -> "Debug using KF. " + new run with additional text ---
$rng = $d.Content
$found = $rng.Find.Execute("Debug using KF", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "Debug using KF. "
    $rng.Collapse(0)
    $rng.InsertAfter("Notice two things to investigate: why in the first step the optimal solution is not the initial solution (in fact, a lot of sensor placement seem to have similar result, as long as the target is in FOV. Check what gradient is used). Also think about a simple case to test the algorithm on.")
    # Force the newly inserted text into its own run (same visual color as original)
    $rng.Font.Color = 10
    $rng.Font.Color = 0xA0000
}

# --- Step 2: remove the whole paragraph "add Kref into the optimization variable" ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("add Kref into the optimization variable", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $p2 = $rng2.Paragraphs.Item(1)
    $wholeRange = $d.Range($p2.Range.Start, $p2.Range.End + 1)
    $wholeRange.Delete()
}
